$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
# row 9 (hunk 0, @@ -1076,22 +1076,22 @@)
$ws.Range("H9").Value = 1603451.4
$ws.Range("I9").Value = 4282
$ws.Range("K9").Value = 4282
$ws.Range("M9").Value = -4113

# row 19 (hunk 1, @@ -1566,25 +1566,25 @@)
$ws.Range("H19").Value = 5501.6665
$ws.Range("J19").Value = 3996
$ws.Range("L19").Value = 3996
$ws.Range("N19").Value = -4346

# row 32 (hunk 2, @@ -2212,22 +2212,22 @@)
$ws.Range("H32").Value = 3077.25
$ws.Range("I32").Value = 2700.3333
$ws.Range("K32").Value = 2700.3333
$ws.Range("M32").Value = -2374.3333

# row 113 (hunk 3, @@ -6259,25 +6259,22 @@)
$ws.Range("H113").Value = 2998.25
$ws.Range("I113").Value = 2998.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2998.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 255.75
$ws.Range("N113").ClearContents()

$ws = $wb.Sheets.Item("ARM")
# row 12 (hunk 4, @@ -8288,25 +8285,25 @@)
$ws.Range("H12").Value = 2310.4
$ws.Range("I12").Value = 1016.3333
$ws.Range("J12").Value = 4251.5
$ws.Range("K12").Value = 1016.3333
$ws.Range("L12").Value = 4251.5
$ws.Range("M12").Value = -843.3333
$ws.Range("N12").Value = -4597.5

# row 16 (hunk 5, @@ -8490,25 +8487,22 @@)
$ws.Range("H16").Value = 198
$ws.Range("I16").Value = 198
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 198
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 89
$ws.Range("N16").ClearContents()

# row 21 (hunk 6, @@ -8732,22 +8726,25 @@)
$ws.Range("H21").Value = 20937.625
$ws.Range("I21").Value = 483.33334
$ws.Range("J21").Value = 33210.2
$ws.Range("K21").Value = 483.33334
$ws.Range("L21").Value = 33210.2
$ws.Range("N21").Value = -33958.2
$ws.Range("M21").Value = -109.33334

# row 30 (hunk 7, @@ -9173,22 +9170,22 @@)
$ws.Range("H30").Value = 1163.4
$ws.Range("I30").Value = 829.25
$ws.Range("K30").Value = 829.25
$ws.Range("M30").Value = -679.25

# row 35 (hunk 8, @@ -9427,25 +9424,25 @@)
$ws.Range("H35").Value = 6688.231
$ws.Range("I35").Value = 179
$ws.Range("J35").Value = 8641
$ws.Range("K35").Value = 179
$ws.Range("L35").Value = 8641
$ws.Range("M35").Value = 227
$ws.Range("N35").Value = -9453

# row 41 (hunk 9, @@ -9730,26 +9727,23 @@)
$ws.Range("H41").Value = 3616
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

# row 42 (hunk 10, @@ -9782,19 +9776,22 @@)
$ws.Range("H42").Value = 10499.5
$ws.Range("I42").Value = 10499.5
$ws.Range("K42").Value = 10499.5
$ws.Range("M42").Value = -10013.5

# row 45 (hunk 11, @@ -9929,25 +9926,25 @@)
$ws.Range("H45").Value = 3478.182
$ws.Range("I45").Value = 2970.125
$ws.Range("J45").Value = 4833
$ws.Range("K45").Value = 2970.125
$ws.Range("L45").Value = 4833
$ws.Range("M45").Value = -2593.125
$ws.Range("N45").Value = -5587

# row 47 (hunk 12, @@ -10030,19 +10027,22 @@)
$ws.Range("H47").Value = 29041
$ws.Range("J47").Value = 29041
$ws.Range("L47").Value = 29041
$ws.Range("N47").Value = -30491

# row 63 (hunk 13, @@ -10805,22 +10805,22 @@)
$ws.Range("H63").Value = 7504.577
$ws.Range("I63").Value = 1457.8572
$ws.Range("K63").Value = 1457.8572
$ws.Range("M63").Value = -771.8571999999999

# row 66 (hunk 14, @@ -10955,22 +10955,22 @@)
$ws.Range("H66").Value = 7504.577
$ws.Range("I66").Value = 1457.8572
$ws.Range("K66").Value = 7289.286
$ws.Range("M66").Value = -3857.286

# row 110 (hunk 15, @@ -13072,22 +13072,19 @@)
$ws.Range("I110").Value = 128334340
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 128334340
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -128332295
$ws.Range("N110").ClearContents()

$ws = $wb.Sheets.Item("BSM")
# row 20 (hunk 16, @@ -15598,25 +15595,25 @@)
$ws.Range("H20").Value = 3344.0908
$ws.Range("I20").Value = 3072
$ws.Range("J20").Value = 3820.25
$ws.Range("K20").Value = 3072
$ws.Range("L20").Value = 3820.25
$ws.Range("M20").Value = -2825
$ws.Range("N20").Value = -4314.25

# row 86 (hunk 17, @@ -18799,25 +18796,25 @@)
$ws.Range("H86").Value = 215435.58
$ws.Range("I86").Value = 4644.727
$ws.Range("J86").Value = 505273
$ws.Range("K86").Value = 4644.727
$ws.Range("L86").Value = 505273
$ws.Range("M86").Value = -3521.727
$ws.Range("N86").Value = -507519

# row 89 (hunk 18, @@ -18946,25 +18943,25 @@)
$ws.Range("H89").Value = 215435.58
$ws.Range("I89").Value = 4644.727
$ws.Range("J89").Value = 505273
$ws.Range("K89").Value = 23223.635
$ws.Range("L89").Value = 2526365
$ws.Range("M89").Value = -17607.635
$ws.Range("N89").Value = -2537597

# row 97 (hunk 19, @@ -19338,22 +19335,22 @@)
$ws.Range("H97").Value = 26495
$ws.Range("I97").Value = 13674
$ws.Range("K97").Value = 13674
$ws.Range("M97").Value = -12683

# row 132 (hunk 20, @@ -21023,22 +21020,22 @@)
$ws.Range("H132").Value = 99074
$ws.Range("J132").Value = 99074
$ws.Range("L132").Value = 99074
$ws.Range("N132").Value = -109194

# row 134 (hunk 21, @@ -21121,22 +21118,22 @@)
$ws.Range("H134").Value = 5664.5293
$ws.Range("I134").Value = 2709.4546
$ws.Range("K134").Value = 8128.3638
$ws.Range("M134").Value = -5593.3638

$ws = $wb.Sheets.Item("CRP")
# row 16 (hunk 22, @@ -22314,25 +22311,25 @@)
$ws.Range("H16").Value = 4232.724
$ws.Range("I16").Value = 2888
$ws.Range("J16").Value = 4837.85
$ws.Range("K16").Value = 2888
$ws.Range("L16").Value = 4837.85
$ws.Range("M16").Value = -2601
$ws.Range("N16").Value = -5411.85

# row 53 (hunk 23, @@ -24175,22 +24172,22 @@)
$ws.Range("H53").Value = 60424.668
$ws.Range("J53").Value = 60424.668
$ws.Range("L53").Value = 60424.668
$ws.Range("N53").Value = -61638.668

# row 58 (hunk 24, @@ -24423,22 +24420,22 @@)
$ws.Range("H58").Value = 5155.364
$ws.Range("I58").Value = 4756.7144
$ws.Range("K58").Value = 4756.7144
$ws.Range("M58").Value = -4553.7144

# row 62 (hunk 25, @@ -24631,22 +24628,22 @@)
$ws.Range("H62").Value = 29999.666
$ws.Range("J62").Value = 29999.666
$ws.Range("L62").Value = 29999.666
$ws.Range("N62").Value = -31247.666

# row 65 (hunk 26, @@ -24772,22 +24769,22 @@)
$ws.Range("H65").Value = 29999.666
$ws.Range("J65").Value = 29999.666
$ws.Range("L65").Value = 149998.33
$ws.Range("N65").Value = -156238.33

# row 99 (hunk 27, @@ -26432,19 +26429,22 @@)
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

# row 113 (hunk 28, @@ -27109,25 +27109,25 @@)
$ws.Range("H113").Value = 4232.724
$ws.Range("I113").Value = 2888
$ws.Range("J113").Value = 4837.85
$ws.Range("K113").Value = 2888
$ws.Range("L113").Value = 4837.85
$ws.Range("M113").Value = -718
$ws.Range("N113").Value = -9177.85

# row 126 (hunk 29, @@ -27734,19 +27734,22 @@)
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

# row 136 (hunk 30, @@ -28221,22 +28224,22 @@)
$ws.Range("H136").Value = 5155.364
$ws.Range("I136").Value = 4756.7144
$ws.Range("K136").Value = 14270.1432
$ws.Range("M136").Value = -11720.1432

$ws = $wb.Sheets.Item("CUL")
# row 38 (hunk 31, @@ -30442,7 +30445,7 @@)
$ws.Range("H38").Value = 116.818184

# row 55 (hunk 32, @@ -31308,22 +31311,22 @@)
$ws.Range("H55").Value = 1067.3636
$ws.Range("I55").Value = 292.8
$ws.Range("K55").Value = 878.4000000000001
$ws.Range("M55").Value = -701.4000000000001

# row 136 (hunk 33, @@ -35403,22 +35406,22 @@)
$ws.Range("H136").Value = 2820
$ws.Range("I136").Value = 2820
$ws.Range("K136").Value = 8460
$ws.Range("M136").Value = -3360

# row 138 (hunk 34, @@ -35504,22 +35507,22 @@)
$ws.Range("H138").Value = 4169178.2
$ws.Range("I138").Value = 6251580
$ws.Range("K138").Value = 18754740
$ws.Range("M138").Value = -18749600

$ws = $wb.Sheets.Item("GSM")
# row 70 (hunk 35, @@ -39171,22 +39174,22 @@)
$ws.Range("H70").Value = 21028.285
$ws.Range("I70").Value = 20290
$ws.Range("K70").Value = 20290
$ws.Range("M70").Value = -20020

# row 73 (hunk 36, @@ -39318,22 +39321,22 @@)
$ws.Range("H73").Value = 21028.285
$ws.Range("I73").Value = 20290
$ws.Range("K73").Value = 20290
$ws.Range("M73").Value = -19354

# row 100 (hunk 37, @@ -40614,22 +40617,22 @@)
$ws.Range("H100").Value = 21900
$ws.Range("J100").Value = 21900
$ws.Range("L100").Value = 21900
$ws.Range("N100").Value = -24064

# row 128 (hunk 38, @@ -41977,22 +41980,22 @@)
$ws.Range("H128").Value = 86999
$ws.Range("J128").Value = 86999
$ws.Range("L128").Value = 86999
$ws.Range("N128").Value = -96959

# row 132 (hunk 39, @@ -42167,25 +42170,25 @@)
$ws.Range("H132").Value = 4976.5625
$ws.Range("J132").Value = 6660.4116
$ws.Range("L132").Value = 19981.2348
$ws.Range("N132").Value = -25041.2348

$ws = $wb.Sheets.Item("LTW")
# row 43 (hunk 40, @@ -44760,22 +44763,19 @@)
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# row 46 (hunk 41, @@ -44901,25 +44901,25 @@)
$ws.Range("H46").Value = 4991.3145
$ws.Range("I46").Value = 1970.8572
$ws.Range("J46").Value = 5746.4287
$ws.Range("K46").Value = 1970.8572
$ws.Range("L46").Value = 5746.4287
$ws.Range("M46").Value = -1782.8572
$ws.Range("N46").Value = -6122.4287

# row 61 (hunk 42, @@ -45651,25 +45651,25 @@)
$ws.Range("H61").Value = 2540
$ws.Range("I61").Value = 2658.3333
$ws.Range("J61").Value = 1297.5
$ws.Range("K61").Value = 2658.3333
$ws.Range("L61").Value = 1297.5
$ws.Range("M61").Value = -2456.3333
$ws.Range("N61").Value = -1701.5

# row 113 (hunk 43, @@ -48184,25 +48184,25 @@)
$ws.Range("H113").Value = 2540
$ws.Range("I113").Value = 2658.3333
$ws.Range("J113").Value = 1297.5
$ws.Range("K113").Value = 2658.3333
$ws.Range("L113").Value = 1297.5
$ws.Range("M113").Value = -488.3332999999998
$ws.Range("N113").Value = -5637.5

$ws = $wb.Sheets.Item("WVR")
# row 61 (hunk 44, @@ -52575,25 +52575,25 @@)
$ws.Range("H61").Value = 14488.2
$ws.Range("I61").Value = 12353.125
$ws.Range("J61").Value = 23028.5
$ws.Range("K61").Value = 12353.125
$ws.Range("L61").Value = 23028.5
$ws.Range("M61").Value = -12061.125
$ws.Range("N61").Value = -23612.5

# row 108 (hunk 45, @@ -54869,19 +54869,22 @@)
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
